# Update "PERIOD TO EXPIRE" (H) and "LAST UPDATE" (I) columns for rows 3-7
# on the "Training Dashboard" sheet to reflect progress as of 04-Nov-2025:
#   - PERIOD TO EXPIRE decreases by one day
#   - LAST UPDATE changes from 03-Nov-2025 to 04-Nov-2025

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Format the "LAST UPDATE" column as text first so the date-like string
# we write ("04-Nov-2025") is kept as literal text instead of being
# auto-converted into a date serial number.
$lastUpdateRange = $ws.Range("I3:I7")
$lastUpdateRange.NumberFormat = "@"

for ($row = 3; $row -le 7; $row++) {
    $periodCell = $ws.Cells.Item($row, 8)   # column H: PERIOD TO EXPIRE
    $periodCell.Value2 = $periodCell.Value2 - 1

    $lastUpdateCell = $ws.Cells.Item($row, 9)   # column I: LAST UPDATE
    $lastUpdateCell.Value2 = "04-Nov-2025"
}
